$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$c = $ws.Range("D2")
$c.NumberFormat = "@"
$c.Value2 = '27.800.56'
$c.Style = "Normal"
$ws.Range("E2").Value2 = '  +0.02%  '

# Row 3
$c = $ws.Range("D3")
$c.NumberFormat = "@"
$c.Value2 = '1.854.70'
$c.Style = "Normal"
$ws.Range("E3").Value2 = '  -0.59%  '

# Row 4
$ws.Range("E4").Value2 = '  -1.96%  '

# Row 5
$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value2 = '320.85'
$c.Style = "Normal"
$ws.Range("E5").Value2 = '  -1.41%  '

# Row 6
$ws.Range("E6").Value2 = '  -2.14%  '

# Row 7
$c = $ws.Range("D7")
$c.NumberFormat = "@"
$c.Value2 = '0.4322'
$c.Style = "Normal"
$ws.Range("E7").Value2 = '  -2.19%  '

# Row 8
$c = $ws.Range("D8")
$c.NumberFormat = "@"
$c.Value2 = '0.3776'
$c.Style = "Normal"
$ws.Range("E8").Value2 = '  -0.61%  '

# Row 9
$c = $ws.Range("D9")
$c.NumberFormat = "@"
$c.Value2 = '0.07410'
$c.Style = "Normal"
$ws.Range("E9").Value2 = '  -0.68%  '

# Row 10
$c = $ws.Range("D10")
$c.NumberFormat = "@"
$c.Value2 = '0.8868'
$c.Style = "Normal"
$ws.Range("E10").Value2 = '  +0.16%  '

# Row 11
$ws.Range("E11").Value2 = '  -0.29%  '

# Row 12
$c = $ws.Range("D12")
$c.NumberFormat = "@"
$c.Value2 = '1.854.76'
$c.Style = "Normal"
$ws.Range("E12").Value2 = '  -0.93%  '

# Row 13
$c = $ws.Range("D13")
$c.NumberFormat = "@"
$c.Value2 = '6.770'
$c.Style = "Normal"

# Row 14
$c = $ws.Range("D14")
$c.NumberFormat = "@"
$c.Value2 = '5.497'
$c.Style = "Normal"
$ws.Range("E14").Value2 = '  -1.15%  '

# Row 15
$c = $ws.Range("D15")
$c.NumberFormat = "@"
$c.Value2 = '0.07142'
$c.Style = "Normal"
$ws.Range("E15").Value2 = '  -1.24%  '

# Row 16
$c = $ws.Range("D16")
$c.NumberFormat = "@"
$c.Value2 = '88.78'
$c.Style = "Normal"
$ws.Range("E16").Value2 = '  +5.99%  '

# Row 17
$ws.Range("E17").Value2 = '  -2.20%  '

# Row 18
$c = $ws.Range("D18")
$c.NumberFormat = "@"
$c.Value2 = '0.000009039'
$c.Style = "Normal"
$ws.Range("E18").Value2 = '  -1.06%  '

# Row 19
$c = $ws.Range("D19")
$c.NumberFormat = "@"
$c.Value2 = '1.014'
$c.Style = "Normal"
$ws.Range("E19").Value2 = '  -2.06%  '

# Row 20
$c = $ws.Range("D20")
$c.NumberFormat = "@"
$c.Value2 = '15.54'
$c.Style = "Normal"
$ws.Range("E20").Value2 = '  -0.14%  '

# Row 21
$c = $ws.Range("D21")
$c.NumberFormat = "@"
$c.Value2 = '27.851.79'
$c.Style = "Normal"
$ws.Range("E21").Value2 = '  +0.11%  '

# Row 22
$c = $ws.Range("D22")
$c.NumberFormat = "@"
$c.Value2 = '5.285'
$c.Style = "Normal"
$ws.Range("E22").Value2 = '  -0.70%  '

# Row 23
$c = $ws.Range("D23")
$c.NumberFormat = "@"
$c.Value2 = '11.21'
$c.Style = "Normal"
$ws.Range("E23").Value2 = '  -1.85%  '

# Row 24
$c = $ws.Range("D24")
$c.NumberFormat = "@"
$c.Value2 = '2.092.03'
$c.Style = "Normal"
$ws.Range("E24").Value2 = '  -0.15%  '

# Row 25
$c = $ws.Range("D25")
$c.NumberFormat = "@"
$c.Value2 = '2.032'
$c.Style = "Normal"
$ws.Range("E25").Value2 = '  +3.35%  '

# Row 26
$c = $ws.Range("D26")
$c.NumberFormat = "@"
$c.Value2 = '156.43'
$c.Style = "Normal"

# Row 27
$ws.Range("E27").Value2 = '  -1.05%  '

# Row 28
$c = $ws.Range("D28")
$c.NumberFormat = "@"
$c.Value2 = '2.107'
$c.Style = "Normal"
$ws.Range("E28").Value2 = '  +5.53%  '

# Row 29
$c = $ws.Range("D29")
$c.NumberFormat = "@"
$c.Value2 = '5.437'
$c.Style = "Normal"
$ws.Range("E29").Value2 = '  +1.96%  '

# Row 30
$c = $ws.Range("D30")
$c.NumberFormat = "@"
$c.Value2 = '121.97'
$c.Style = "Normal"
$ws.Range("E30").Value2 = '  +3.71%  '

# Row 31
$c = $ws.Range("D31")
$c.NumberFormat = "@"
$c.Value2 = '0.08972'
$c.Style = "Normal"
$ws.Range("E31").Value2 = '  -1.51%  '

# Row 32
$c = $ws.Range("D32")
$c.NumberFormat = "@"
$c.Value2 = '1.244'
$c.Style = "Normal"
$ws.Range("E32").Value2 = '  +2.15%  '

# Row 33
$c = $ws.Range("D33")
$c.NumberFormat = "@"
$c.Value2 = '0.7820'
$c.Style = "Normal"
$ws.Range("E33").Value2 = '  +0.53%  '

# Row 34
$c = $ws.Range("D34")
$c.NumberFormat = "@"
$c.Value2 = '4.587'
$c.Style = "Normal"
$ws.Range("E34").Value2 = '  -0.07%  '

# Row 35
$c = $ws.Range("D35")
$c.NumberFormat = "@"
$c.Value2 = '2.930'
$c.Style = "Normal"
$ws.Range("E35").Value2 = '  -4.01%  '

# Row 36
$ws.Range("E36").Value2 = '  -1.54%  '

# Row 37
$ws.Range("E37").Value2 = '  -2.07%  '

# Row 38
$c = $ws.Range("D38")
$c.NumberFormat = "@"
$c.Value2 = '0.05336'
$c.Style = "Normal"
$ws.Range("E38").Value2 = '  -0.51%  '

# Row 39
$c = $ws.Range("D39")
$c.NumberFormat = "@"
$c.Value2 = '0.01973'
$c.Style = "Normal"
$ws.Range("E39").Value2 = '  -1.16%  '

# Row 40
$c = $ws.Range("D40")
$c.NumberFormat = "@"
$c.Value2 = '2.881'
$c.Style = "Normal"
$ws.Range("E40").Value2 = '  +1.44%  '

# Row 41
$ws.Range("B41").Value2 = 'FraxShare'
$ws.Range("C41").Value2 = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$c = $ws.Range("D41")
$c.NumberFormat = "@"
$c.Value2 = '7.088'
$c.Style = "Normal"
$ws.Range("E41").Value2 = '  +2.63%  '

# Row 42
$ws.Range("B42").Value2 = 'TheSandbox'
$ws.Range("C42").Value2 = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$c = $ws.Range("D42")
$c.NumberFormat = "@"
$c.Value2 = '0.5210'
$c.Style = "Normal"
$ws.Range("E42").Value2 = '  -0.09%  '

# Row 43
$c = $ws.Range("D43")
$c.NumberFormat = "@"
$c.Value2 = '0.1688'
$c.Style = "Normal"
$ws.Range("E43").Value2 = '  -0.49%  '

# Row 44
$c = $ws.Range("D44")
$c.NumberFormat = "@"
$c.Value2 = '8.889'
$c.Style = "Normal"
$ws.Range("E44").Value2 = '  +1.83%  '

# Row 45
$c = $ws.Range("D45")
$c.NumberFormat = "@"
$c.Value2 = '111.21'
$c.Style = "Normal"
$ws.Range("E45").Value2 = '  +1.27%  '

# Row 46
$c = $ws.Range("D46")
$c.NumberFormat = "@"
$c.Value2 = '10.72'
$c.Style = "Normal"
$ws.Range("E46").Value2 = '  +0.68%  '

# Row 47
$c = $ws.Range("D47")
$c.NumberFormat = "@"
$c.Value2 = '1.717'
$c.Style = "Normal"
$ws.Range("E47").Value2 = '  -0.53%  '

# Row 48
$c = $ws.Range("D48")
$c.NumberFormat = "@"
$c.Value2 = '0.4755'
$c.Style = "Normal"
$ws.Range("E48").Value2 = '  +0.86%  '

# Row 49
$c = $ws.Range("D49")
$c.NumberFormat = "@"
$c.Value2 = '0.06523'
$c.Style = "Normal"
$ws.Range("E49").Value2 = '  +1.31%  '

# Row 50
$c = $ws.Range("D50")
$c.NumberFormat = "@"
$c.Value2 = '1.015'
$c.Style = "Normal"
$ws.Range("E50").Value2 = '  -2.31%  '

# Row 51
$c = $ws.Range("D51")
$c.NumberFormat = "@"
$c.Value2 = '1.903'
$c.Style = "Normal"
$ws.Range("E51").Value2 = '  +0.78%  '
